$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = -0.1824983152294041
$ws.Range("D2").Value = 0.8562755077668034

$ws.Range("C3").Value = 0.529181002213955
$ws.Range("D3").Value = 0.6001181404730382

$ws.Range("C4").Value = -0.03453432281882979
$ws.Range("D4").Value = 0.9726530265765045

$ws.Range("C5").Value = -1.523546762671858
$ws.Range("D5").Value = 0.136870152946527

$ws.Range("C6").Value = 0.576336479316639
$ws.Range("D6").Value = 0.5681838106766759

$ws.Range("C7").Value = 0.09597620465980004
$ws.Range("D7").Value = 0.9241029518723634

$ws.Range("C8").Value = -1.684867181147551
$ws.Range("D8").Value = 0.1011733945798423

$ws.Range("C9").Value = -0.3441934839789854
$ws.Range("D9").Value = 0.7328205689259306

$ws.Range("C10").Value = -1.943634537787615
$ws.Range("D10").Value = 0.0602540882740632

$ws.Range("C11").Value = -1.429762543870582
$ws.Range("D11").Value = 0.161910506926948
